$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Triggers")

# Insert a new "Classifier" column before the existing "Skill mean" column (E).
$ws.Columns("E:E").Insert()

# The inserted column picks up the formatting of its neighbour; strip that so
# the new cells come in unstyled (matching a freshly typed-in column).
$ws.Range("E1:E8").ClearFormats()

# Header
$ws.Cells.Item(1, 5).Value2 = "Classifier"

# Per-row classifier labels (first-use order here matches the shared-string
# table order of the saved workbook: Classifier, Gaussian, L1T, Dummy).
$ws.Cells.Item(5, 5).Value2 = "Gaussian"
$ws.Cells.Item(6, 5).Value2 = "L1T"
$ws.Cells.Item(2, 5).Value2 = "Dummy"
$ws.Cells.Item(3, 5).Value2 = "Dummy"
$ws.Cells.Item(4, 5).Value2 = "Dummy"
$ws.Cells.Item(7, 5).Value2 = "Gaussian"
$ws.Cells.Item(8, 5).Value2 = "Dummy"

# Row 6 (L1T) had its old "Skill mean"/"Skill variance" values (now columns F/G)
# reset to 0 as part of this working revision.
$ws.Cells.Item(6, 6).Value2 = 0
$ws.Cells.Item(6, 7).Value2 = 0

# Update the active selection to match the saved workbook state.
$ws.Range("G11").Select()
